# Update "想去人数" (want-to-go count) values in column F across all four
# sheets of the workbook, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 3418
$ws1.Range("F5").Value  = 229
$ws1.Range("F6").Value  = 4920
$ws1.Range("F7").Value  = 489
$ws1.Range("F8").Value  = 318
$ws1.Range("F11").Value = 289
$ws1.Range("F12").Value = 57
$ws1.Range("F15").Value = 296
$ws1.Range("F20").Value = 351
$ws1.Range("F21").Value = 4801
$ws1.Range("F22").Value = 32
$ws1.Range("F25").Value = 5932
$ws1.Range("F27").Value = 7
$ws1.Range("F28").Value = 3204
$ws1.Range("F29").Value = 279
$ws1.Range("F30").Value = 686
$ws1.Range("F31").Value = 4428
$ws1.Range("F33").Value = 102
$ws1.Range("F35").Value = 905
$ws1.Range("F37").Value = 16
$ws1.Range("F39").Value = 819
$ws1.Range("F40").Value = 899
$ws1.Range("F41").Value = 9

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 39
$ws2.Range("F6").Value = 51

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 1098
$ws3.Range("F4").Value = 43

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 1098
$ws4.Range("F5").Value  = 43
$ws4.Range("F8").Value  = 3419
$ws4.Range("F9").Value  = 229
$ws4.Range("F10").Value = 4920
$ws4.Range("F11").Value = 489
$ws4.Range("F12").Value = 318
$ws4.Range("F15").Value = 289
$ws4.Range("F16").Value = 57
$ws4.Range("F19").Value = 296
$ws4.Range("F21").Value = 39
$ws4.Range("F25").Value = 351
$ws4.Range("F26").Value = 4801
$ws4.Range("F27").Value = 32
$ws4.Range("F30").Value = 5932
$ws4.Range("F32").Value = 7
$ws4.Range("F33").Value = 3204
$ws4.Range("F34").Value = 279
$ws4.Range("F35").Value = 686
$ws4.Range("F36").Value = 4428
$ws4.Range("F39").Value = 102
$ws4.Range("F40").Value = 905
$ws4.Range("F42").Value = 16
$ws4.Range("F44").Value = 819
$ws4.Range("F45").Value = 899
$ws4.Range("F47").Value = 9
$ws4.Range("F49").Value = 51
